# edit.ps1 - PowerPoint COM-interop script (PowerShell-style)
# Applies the authored content changes:
#   1. Slide 11 ("Next Steps"), Content Placeholder 2, 3rd bullet:
#        "In " + "WG adoption (SPRING WG) queue"
#        -> "Request SPRING " + "WG adoption"
#   2. Slide 6, shape "Rectangle 3": horizontal offset 3681950 EMU -> 3657600 EMU
#      (i.e. Left = 288.0 points)

$p = $ppt.ActivePresentation

# --- 1) Slide 11: update the "WG adoption" bullet text -------------------
$slide11 = $p.Slides.Item(11)
$contentShape = $slide11.Shapes.Item(2)   # "Content Placeholder 2"
$textRange = $contentShape.TextFrame.TextRange

$bulletPara = $textRange.Paragraphs(3, 1)
$bulletPara.Runs(1, 1).Text = "Request SPRING "
$bulletPara.Runs(2, 1).Text = "WG adoption"

# --- 2) Slide 6: reposition "Rectangle 3" ---------------------------------
$slide6 = $p.Slides.Item(6)
$rect3 = $slide6.Shapes.Item(5)           # "Rectangle 3"
$rect3.Left = 288.0
